# Add a "Problem Sheet" note textbox to the "Bibliography Management"
# slide (sldId 322 / creationId 687175056), mirroring the author's edit.

$p = $ppt.ActivePresentation

# Locate the target slide by its persistent SlideID (322) rather than a
# hard-coded index, in case slide ordering ever changes.
$slide = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    if ($p.Slides.Item($i).SlideID -eq 322) {
        $slide = $p.Slides.Item($i)
        break
    }
}
if ($slide -eq $null) {
    $slide = $p.Slides.Item(24)
}

# Position/size taken straight from the target OOXML (EMU -> points,
# 914400 EMU per inch / 12700 EMU per point).
$left   = 7357241 / 12700
$top    = 6498021 / 12700
$width  = 1697422 / 12700
$height = 369332 / 12700

$tb = $slide.Shapes.AddTextbox(1, $left, $top, $width, $height)

# Match the original shape formatting: auto-fit textbox, no fill, and
# accent2-colored text reading "Problem Sheet".
$tb.Fill.Visible = $false
$tb.TextFrame.WordWrap = -1
$tb.TextFrame.AutoSize = 1

$tr = $tb.TextFrame.TextRange
$tr.Text = "Problem Sheet"
$tr.LanguageID = "en-GB"
$tr.Font.Color.ObjectThemeColor = 6
